$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append ".age_trait" to every ".deja" header label in row 1, columns B through P.
$headerCells = @("B1","C1","D1","E1","F1","G1","H1","I1","J1","K1","L1","M1","N1","O1","P1")

foreach ($addr in $headerCells) {
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2 + ".age_trait"
}
